$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking price cells so Excel does not
# auto-convert them to numbers (these are text cells in the source data).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = '26.687.42'
$ws.Range("E2").Value = '  +1.44%  '
$ws.Range("D3").Value = '1.631.40'
$ws.Range("E3").Value = '  +1.37%  '
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").Value = '213.65'
$ws.Range("E5").Value = '  +0.44%  '
$ws.Range("D7").Value = '0.490'
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("E8").Value = '  +0.75%  '
$ws.Range("D9").Value = '0.0620'
$ws.Range("E9").Value = '  +0.75%  '
$ws.Range("D10").Value = '19.03'
$ws.Range("E10").Value = '  +4.23%  '
$ws.Range("D11").Value = '0.0834'
$ws.Range("E11").Value = '  +2.40%  '
$ws.Range("D12").Value = '1.856.87'
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("D13").Value = '1.627.02'
$ws.Range("E13").Value = '  +1.16%  '
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("E15").Value = '  +1.97%  '
$ws.Range("D16").Value = '26.643.64'
$ws.Range("E16").Value = '  +1.36%  '
$ws.Range("D17").Value = '63.01'
$ws.Range("E17").Value = '  +2.29%  '
$ws.Range("D18").Value = '0.0₃0733'
$ws.Range("E18").Value = '  +0.49%  '
$ws.Range("B19").Value = 'BitcoinCash'
$ws.Range("C19").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D19").Value = '209.30'
$ws.Range("E19").Value = '  +2.80%  '
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").Value = '1.00'
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").Value = '4.33'
$ws.Range("E21").Value = '  +0.72%  '
$ws.Range("E22").Value = '  +1.05%  '
$ws.Range("D23").Value = '6.09'
$ws.Range("E23").Value = '  +1.41%  '
$ws.Range("E24").Value = '  -2.01%  '
$ws.Range("D25").Value = '145.82'
$ws.Range("E25").Value = '  +1.04%  '
$ws.Range("E26").Value = '  +0.02%  '
$ws.Range("E27").Value = '  -1.42%  '
$ws.Range("D28").Value = '15.37'
$ws.Range("E28").Value = '  +0.85%  '
$ws.Range("E29").Value = '  +1.16%  '
$ws.Range("E30").Value = '  +6.37%  '
$ws.Range("D31").Value = '1.17'
$ws.Range("E31").Value = '  +0.22%  '
$ws.Range("D32").Value = '3.23'
$ws.Range("E32").Value = '  +0.96%  '
$ws.Range("D33").Value = '2.96'
$ws.Range("E33").Value = '  +0.84%  '
$ws.Range("D34").Value = '1.51'
$ws.Range("E34").Value = '  +1.28%  '
$ws.Range("E35").Value = '  -0.28%  '
$ws.Range("D36").Value = '1.164.77'
$ws.Range("E36").Value = '  +0.85%  '
$ws.Range("D37").Value = '0.0166'
$ws.Range("E37").Value = '  +0.90%  '
$ws.Range("D38").Value = '0.813'
$ws.Range("E38").Value = '  +2.11%  '
$ws.Range("E39").Value = '  +0.11%  '
$ws.Range("E40").Value = '  -0.32%  '
$ws.Range("D41").Value = '0.502'
$ws.Range("E41").Value = '  +0.31%  '
$ws.Range("D42").Value = '5.41'
$ws.Range("E42").Value = '  +3.00%  '
$ws.Range("D43").Value = '0.786'
$ws.Range("E43").Value = '  +0.21%  '
$ws.Range("D44").Value = '1.767.44'
$ws.Range("E44").Value = '  +1.49%  '
$ws.Range("D45").Value = '92.20'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("E46").Value = '  +1.68%  '
$ws.Range("D47").Value = '54.63'
$ws.Range("E47").Value = '  +0.61%  '
$ws.Range("D48").Value = '0.0512'
$ws.Range("E48").Value = '  +1.14%  '
$ws.Range("D49").Value = '7.59'
$ws.Range("E49").Value = '  +4.64%  '
$ws.Range("E50").Value = '  +0.78%  '
$ws.Range("E51").Value = '  +0.17%  '
